$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3559676545.5986443
$ws.Range("C3").Value = 3565451560.217553
$ws.Range("C4").Value = 3577001625.702943
$ws.Range("C5").Value = 3590861768.0813828
